# Auto-generated edit script applying the cell-value changes described in the diff.
# Each FFXIV leve-profit sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) gets a handful
# of market-price driven recalculations: most cells are simple value updates, a few
# previously-blank profit cells now carry a value, and a few previously-populated cells
# are cleared back to blank (to mirror the OOXML no longer emitting that <c> element).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5300
$ws.Range("I86").Value = 1900
$ws.Range("J86").Value = 7000
$ws.Range("K86").Value = 1900
$ws.Range("L86").Value = 7000
$ws.Range("M86").Value = -777
$ws.Range("N86").Value = -9246
$ws.Range("H87").Value = 74652
$ws.Range("J87").Value = 74652
$ws.Range("L87").Value = 74652
$ws.Range("N87").Value = -77148
$ws.Range("H89").Value = 5300
$ws.Range("I89").Value = 1900
$ws.Range("J89").Value = 7000
$ws.Range("K89").Value = 9500
$ws.Range("L89").Value = 35000
$ws.Range("M89").Value = -3884
$ws.Range("N89").Value = -46232
$ws.Range("H90").Value = 74652
$ws.Range("J90").Value = 74652
$ws.Range("L90").Value = 223956
$ws.Range("N90").Value = -236436
$ws.Range("H92").Value = 1355.05
$ws.Range("J92").Value = 1313.5714
$ws.Range("L92").Value = 1313.5714
$ws.Range("N92").Value = -3809.5714
$ws.Range("H113").Value = 3036
$ws.Range("I113").Value = 2398
$ws.Range("J113").Value = 4950
$ws.Range("K113").Value = 2398
$ws.Range("L113").Value = 4950
$ws.Range("M113").Value = 856
$ws.Range("N113").Value = -11458

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1699
$ws.Range("I2").Value = 1699
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1699
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1586
$ws.Range("N2").ClearContents()
$ws.Range("H4").Value = 88
$ws.Range("J4").Value = 200
$ws.Range("L4").Value = 200
$ws.Range("N4").Value = -432
$ws.Range("H74").Value = 9262.4
$ws.Range("I74").Value = 7291.6
$ws.Range("K74").Value = 7291.6
$ws.Range("M74").Value = -6417.6
$ws.Range("H77").Value = 9262.4
$ws.Range("I77").Value = 7291.6
$ws.Range("K77").Value = 36458
$ws.Range("M77").Value = -32090
$ws.Range("H116").Value = 1699
$ws.Range("I116").Value = 1699
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1699
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 595
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1699
$ws.Range("I3").Value = 1699
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1699
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1585
$ws.Range("N3").ClearContents()
$ws.Range("H64").Value = 3283.5715
$ws.Range("I64").Value = 1962.6666
$ws.Range("J64").Value = 4274.25
$ws.Range("K64").Value = 1962.6666
$ws.Range("L64").Value = 4274.25
$ws.Range("M64").Value = -1737.6666
$ws.Range("N64").Value = -4724.25
$ws.Range("H67").Value = 3283.5715
$ws.Range("I67").Value = 1962.6666
$ws.Range("J67").Value = 4274.25
$ws.Range("K67").Value = 1962.6666
$ws.Range("L67").Value = 4274.25
$ws.Range("M67").Value = -1182.6666
$ws.Range("N67").Value = -5834.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 227.5
$ws.Range("I7").Value = 156.75
$ws.Range("J7").Value = 298.25
$ws.Range("K7").Value = 156.75
$ws.Range("L7").Value = 298.25
$ws.Range("M7").Value = -43.75
$ws.Range("N7").Value = -524.25
$ws.Range("H74").Value = 70209.336
$ws.Range("J74").Value = 70209.336
$ws.Range("L74").Value = 70209.336
$ws.Range("N74").Value = -71957.336
$ws.Range("H77").Value = 70209.336
$ws.Range("J77").Value = 70209.336
$ws.Range("L77").Value = 210628.008
$ws.Range("N77").Value = -219364.008
$ws.Range("H105").Value = 2937.8
$ws.Range("I105").Value = 2937.8
$ws.Range("K105").Value = 2937.8
$ws.Range("M105").Value = -1190.8
$ws.Range("H107").Value = 1416.909
$ws.Range("I107").Value = 965.3333
$ws.Range("K107").Value = 965.3333
$ws.Range("M107").Value = 954.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 500
$ws.Range("J92").Value = 500
$ws.Range("L92").Value = 1500
$ws.Range("N92").Value = -3996
$ws.Range("H116").Value = 4785.7144
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H117").Value = 2998
$ws.Range("I117").Value = 2998
$ws.Range("K117").Value = 8994
$ws.Range("M117").Value = -5552
$ws.Range("H131").Value = 4333.3335
$ws.Range("J131").Value = 4333.3335
$ws.Range("L131").Value = 13000.0005
$ws.Range("N131").Value = -23080.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 294.8
$ws.Range("I2").Value = 308.16666
$ws.Range("J2").Value = 274.75
$ws.Range("K2").Value = 308.16666
$ws.Range("L2").Value = 274.75
$ws.Range("M2").Value = -195.16666
$ws.Range("N2").Value = -500.75
$ws.Range("H80").Value = 2776.8
$ws.Range("J80").Value = 3296.6667
$ws.Range("L80").Value = 3296.6667
$ws.Range("N80").Value = -5292.6667
$ws.Range("H83").Value = 2776.8
$ws.Range("J83").Value = 3296.6667
$ws.Range("L83").Value = 16483.3335
$ws.Range("N83").Value = -26467.3335
$ws.Range("H113").Value = 1984.3846
$ws.Range("I113").Value = 1560.7778
$ws.Range("J113").Value = 2937.5
$ws.Range("K113").Value = 1560.7778
$ws.Range("L113").Value = 2937.5
$ws.Range("M113").Value = 609.2221999999999
$ws.Range("N113").Value = -7277.5
$ws.Range("H132").Value = 7498.5
$ws.Range("I132").Value = 4997.5
$ws.Range("K132").Value = 14992.5
$ws.Range("M132").Value = -12462.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1775.6
$ws.Range("I61").Value = 1701.3334
$ws.Range("K61").Value = 1701.3334
$ws.Range("M61").Value = -1499.3334
$ws.Range("H68").Value = 2183.1667
$ws.Range("I68").Value = 2019.8
$ws.Range("K68").Value = 2019.8
$ws.Range("M68").Value = -1270.8
$ws.Range("H71").Value = 2183.1667
$ws.Range("I71").Value = 2019.8
$ws.Range("K71").Value = 10099
$ws.Range("M71").Value = -6355
$ws.Range("H103").Value = 23146.25
$ws.Range("J103").Value = 23146.25
$ws.Range("L103").Value = 23146.25
$ws.Range("N103").Value = -25490.25
$ws.Range("H113").Value = 1775.6
$ws.Range("I113").Value = 1701.3334
$ws.Range("K113").Value = 1701.3334
$ws.Range("M113").Value = 468.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 8000
$ws.Range("I38").Value = 8000
$ws.Range("K38").Value = 8000
$ws.Range("M38").Value = -7527
$ws.Range("H81").Value = 2799
$ws.Range("I81").Value = 2948.6667
$ws.Range("K81").Value = 5897.3334
$ws.Range("M81").Value = -4836.3334
$ws.Range("H84").Value = 2799
$ws.Range("I84").Value = 2948.6667
$ws.Range("K84").Value = 29486.667
$ws.Range("M84").Value = -24182.667
$ws.Range("H113").Value = 858.55
$ws.Range("I113").Value = 932.75
$ws.Range("J113").Value = 747.25
$ws.Range("K113").Value = 2798.25
$ws.Range("L113").Value = 2241.75
$ws.Range("M113").Value = -628.25
$ws.Range("N113").Value = -6581.75
$ws.Range("H126").Value = 1694.6842
$ws.Range("I126").Value = 1600
$ws.Range("K126").Value = 4800
$ws.Range("M126").Value = -2330
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
